# 4/7 minor transition bugs
# - Rename the first sheet ("Comp 380") to "2"
# - Update the Section Number (B1) and Course Color RGB values (D1:F1)
#   on that sheet. The cells hold numeric-looking text (shared strings),
#   so each target cell is explicitly formatted as Text ("@") before the
#   value is typed in, which is what keeps Excel from auto-converting the
#   entry into a real number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Comp 380" -> "2"
$ws.Name = "2"

# Force text storage for the numeric-looking values so they stay strings
$ws.Range("B1").NumberFormat = "@"
$ws.Range("D1").NumberFormat = "@"
$ws.Range("E1").NumberFormat = "@"
$ws.Range("F1").NumberFormat = "@"

# Section Number: 12344 -> 2
$ws.Range("B1").Value = "2"

# Course Color RGB: 204,102,51 -> 0,0,128
$ws.Range("D1").Value = "0"
$ws.Range("E1").Value = "0"
$ws.Range("F1").Value = "128"
